$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "centroid"
$ws.Range("B6").Value = 0.165
$ws.Range("C6").Value = 0.889
$ws.Range("D6").Value = 0.888
$ws.Range("E6").Value = 0.886

$ws.Range("A7").Value = "SVM"
$ws.Range("B7").Value = 0.979
$ws.Range("C7").Value = 0.978
$ws.Range("D7").Value = 0.965
$ws.Range("E7").Value = 0.96

$ws.Range("E12").Select()
